# Generate Report for Handoff
#
# This script moves the localization-status report from the "handed back"
# state to the "ready for handoff" state:
#   - Status cells that said "Handed back: in sync with en-US" now say
#     "Ready for handoff"
#   - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#     timestamps are bumped forward to reflect the new handoff
#   - The (now shorter) status column is narrowed to fit the new text

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-28 04:58:23"

# Narrow the zh-cn / de-de status columns (E, F) to match the shorter text
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-28 04:58:19"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
